$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns E and L (vm_pu values) for rows 2-25 with newly computed
# results from the 380 kV case run.

$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("L20").Value = 0.9929783193494216
$ws.Range("E21").Value = 0.9882828385668253
$ws.Range("L21").Value = 0.9920501090198105
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("L25").Value = 0.9944092447426414
